$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.384.55'
$ws.Range("E2").Value = '  -0.35%  '

$ws.Range("D3").Value = '3.527.45'
$ws.Range("E3").Value = '  -0.94%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.84'
$ws.Range("E5").Value = '  +0.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.73'
$ws.Range("E6").Value = '  -1.56%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.617'
$ws.Range("E7").Value = '  -2.03%  '

$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '3.520.37'
$ws.Range("E8").Value = '  -1.16%  '

$ws.Range("E9").Value = '  +0.15%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.186'
$ws.Range("E10").Value = '  +6.05%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.640'
$ws.Range("E11").Value = '  -2.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.98'
$ws.Range("E12").Value = '  -2.73%  '

$ws.Range("E13").Value = '  +1.82%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.50'
$ws.Range("E14").Value = '  -0.94%  '

$ws.Range("D15").Value = '4.096.02'
$ws.Range("E15").Value = '  -0.82%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.36'
$ws.Range("E16").Value = '  -2.27%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.524.57'
$ws.Range("E17").Value = '  -1.06%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '69.222.46'
$ws.Range("E18").Value = '  -0.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.54'
$ws.Range("E19").Value = '  +0.57%  '

$ws.Range("E20").Value = '  -0.83%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '538.24'
$ws.Range("E21").Value = '  +14.65%  '

$ws.Range("E22").Value = '  +0.19%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '20.80'
$ws.Range("E23").Value = '  +8.94%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.02'
$ws.Range("E24").Value = '  -0.96%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.44'
$ws.Range("E25").Value = '  +3.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '94.74'
$ws.Range("E26").Value = '  +7.62%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.02'
$ws.Range("E27").Value = '  +1.17%  '

$ws.Range("E28").Value = '  -3.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.18'
$ws.Range("E29").Value = '  -0.79%  '

$ws.Range("E30").Value = '  -1.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.28'
$ws.Range("E31").Value = '  -3.80%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.71'
$ws.Range("E32").Value = '  +5.95%  '

$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '64.36'
$ws.Range("E33").Value = '  -1.65%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.114'
$ws.Range("E34").Value = '  -3.46%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '574.66'
$ws.Range("E35").Value = '  +0.32%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.11'
$ws.Range("E36").Value = '  +10.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '38.31'
$ws.Range("E37").Value = '  +0.04%  '

$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  +0.07%  '

$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.400'
$ws.Range("E39").Value = '  +1.94%  '

$ws.Range("D40").Value = '0.0₃0767'
$ws.Range("E40").Value = '  -3.40%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.14'
$ws.Range("E41").Value = '  +2.50%  '

$ws.Range("E42").Value = '  -2.44%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.35'
$ws.Range("E43").Value = '  -3.53%  '

$ws.Range("E44").Value = '  +6.46%  '

$ws.Range("E45").Value = '  -3.75%  '

$ws.Range("D46").Value = '3.199.06'
$ws.Range("E46").Value = '  -0.05%  '

$ws.Range("E47").Value = '  +1.19%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.23'
$ws.Range("E48").Value = '  -2.08%  '

$ws.Range("E49").Value = '  -1.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.996'
$ws.Range("E50").Value = '  -0.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '136.33'
$ws.Range("E51").Value = '  -0.16%  '
